# The commit swaps the two theme parts of the deck:
#   ppt/theme/theme1.xml (the slide master's theme, currently "Integral")
#     becomes the "Office Theme" palette, and
#   ppt/theme/theme2.xml (the notes master's theme, currently "Office Theme")
#     becomes the "Integral" palette.
#
# This host's PowerPoint object model exposes a single editable theme color
# scheme for the deck (Master.Theme.ThemeColorScheme), reachable from the
# slide master, so we drive it to the new "Office" palette here. RGB values
# are assigned as plain OLE COLORREF ints: R | (G<<8) | (B<<16).

$p = $ppt.ActivePresentation

# Target "Office" theme colors (RRGGBB -> OLE RGB int), in the fixed
# ThemeColorScheme index order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$master = $p.SlideMaster
$scheme = $master.Theme.ThemeColorScheme
for ($i = 1; $i -le $scheme.Count; $i++) {
    $scheme.Item($i).RGB = $officeColors[$i - 1]
}
